$wb = $excel.ActiveWorkbook

# ---- Sheet index 1 ----
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 7741
$ws.Range("F3").Value = 7562
$ws.Range("F10").Value = 144
$ws.Range("F11").Value = 219

$ws.Range("C12").Value = "合肥·Look Look动漫嘉年华"
$ws.Range("D12").Value = "新站区东方大道288号 少荃体育中心"
$ws.Range("E12").Value = "2024.05.01 10:00-05.01 17:30"
$ws.Range("F12").Value = 681
$ws.Range("G12").Value = 58
$ws.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=82311"
$ws.Range("I12").Value = "//i2.hdslb.com/bfs/openplatform/202403/jbUNtkAQ1709619599897.png"

$ws.Range("C13").Value = "合肥·第十三届合肥次元之门动漫游戏博览会-多多poi&Mace专场"
$ws.Range("D13").Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws.Range("E13").Value = "2024.05.01 10:00-05.01 17:00"
$ws.Range("F13").Value = 105
$ws.Range("G13").Value = 238
$ws.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=83039"
$ws.Range("I13").Value = "//i0.hdslb.com/bfs/openplatform/202403/8ZIG73sO1710901741120.jpeg"

$ws.Range("F14").Value = 1136
$ws.Range("F16").Value = 40
$ws.Range("F17").Value = 6
$ws.Range("F19").Value = 92

# ---- Sheet index 4 ----
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 7741
$ws.Range("F3").Value = 7563
$ws.Range("F10").Value = 144
$ws.Range("F11").Value = 219

$ws.Range("C12").Value = "合肥·Look Look动漫嘉年华"
$ws.Range("D12").Value = "新站区东方大道288号 少荃体育中心"
$ws.Range("E12").Value = "2024.05.01 10:00-05.01 17:30"
$ws.Range("F12").Value = 681
$ws.Range("G12").Value = 58
$ws.Range("H12").Value = "https://show.bilibili.com/platform/detail.html?id=82311"
$ws.Range("I12").Value = "//i2.hdslb.com/bfs/openplatform/202403/jbUNtkAQ1709619599897.png"

$ws.Range("C13").Value = "合肥·第十三届合肥次元之门动漫游戏博览会-多多poi&Mace专场"
$ws.Range("D13").Value = "南京路与庐州大道交汇处 合肥滨湖国际会展中心"
$ws.Range("E13").Value = "2024.05.01 10:00-05.01 17:00"
$ws.Range("F13").Value = 105
$ws.Range("G13").Value = 238
$ws.Range("H13").Value = "https://show.bilibili.com/platform/detail.html?id=83039"
$ws.Range("I13").Value = "//i0.hdslb.com/bfs/openplatform/202403/8ZIG73sO1710901741120.jpeg"

$ws.Range("F14").Value = 1136
$ws.Range("F16").Value = 40
$ws.Range("F17").Value = 6
$ws.Range("F19").Value = 92

